# Applies the diff: rows 2, 3 and 4 have their data cyclically rotated
# (new row2 <= old row3, new row3 <= old row4, new row4 <= old row2),
# expressed here as the concrete set of per-cell value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 changes ---
$ws.Range("A2").Value = 111697304
$ws.Range("Q2").Value = 373090.8741807578
$ws.Range("R2").Value = 6865424.499624529
$ws.Range("Z2").Value = "19:00"
$ws.Range("AB2").Value = "19:00"
$ws.Range("AI2").Value = "Luckig tallskog. K-skog"
$ws.Range("AJ2").Value = ""
$ws.Range("AK2").Value = ""
$ws.Range("AO2").Value = ""

# --- Row 3 changes ---
$ws.Range("A3").Value = 111697636
$ws.Range("B3").Value = 88489
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 1962
$ws.Range("F3").Value = "Vaddporing"
$ws.Range("G3").Value = "Anomoporia kamtschatica"
$ws.Range("H3").Value = "(Parmasto) Bondartseva"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("Q3").Value = 373112.5181173298
$ws.Range("R3").Value = 6865358.590016441
$ws.Range("AC3").Value = "Växer under rötad gammal silverved"
$ws.Range("AI3").Value = "Kontinuitetsskog. Tallskog"
$ws.Range("AJ3").Value = "tall"
$ws.Range("AK3").Value = "Pinus sylvestris"
$ws.Range("AO3").Value = "Pinus sylvestris"

# --- Row 4 changes ---
$ws.Range("A4").Value = 111697236
$ws.Range("B4").Value = 8377
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 106545
$ws.Range("F4").Value = "Mindre märgborre"
$ws.Range("G4").Value = "Tomicus minor"
$ws.Range("H4").Value = "(Hartig, 1834)"
$ws.Range("J4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "färska gnagspår"
$ws.Range("Q4").Value = 373121.3523494597
$ws.Range("R4").Value = 6865443.651501717
$ws.Range("Z4").Value = "00:00"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AC4").Value = ""
$ws.Range("AI4").Value = "Tallskog. Kontinuitetsskog"
